$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '92.431.32'
$ws.Range("E2").Value = '  +5.60%  '
$ws.Range("D3").Value = '3.288.34'
$ws.Range("E3").Value = '  +0.59%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.34%  '
$ws.Range("D5").Value = '''216.45'
$ws.Range("E5").Value = '  +1.96%  '
$ws.Range("D6").Value = '''634.56'
$ws.Range("E6").Value = '  +1.42%  '
$ws.Range("D7").Value = '''0.411'
$ws.Range("E7").Value = '  +12.32%  '
$ws.Range("D8").Value = '''0.718'
$ws.Range("E8").Value = '  +5.18%  '
$ws.Range("D9").Value = '''0.999'
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").Value = '3.280.32'
$ws.Range("E10").Value = '  +0.39%  '
$ws.Range("E11").Value = '  +2.96%  '
$ws.Range("D12").Value = '''0.0000264'
$ws.Range("E12").Value = '  +3.48%  '
$ws.Range("E13").Value = '  +0.80%  '
$ws.Range("D14").Value = '''34.24'
$ws.Range("E14").Value = '  +1.90%  '
$ws.Range("D15").Value = '3.892.57'
$ws.Range("E15").Value = '  +0.90%  '
$ws.Range("D16").Value = '92.062.47'
$ws.Range("E16").Value = '  +5.86%  '
$ws.Range("E17").Value = '  +1.65%  '
$ws.Range("D18").Value = '3.282.14'
$ws.Range("E18").Value = '  +1.21%  '
$ws.Range("D19").Value = '''3.33'
$ws.Range("E19").Value = '  +7.19%  '
$ws.Range("D20").Value = '''14.09'
$ws.Range("E20").Value = '  +0.67%  '
$ws.Range("E21").Value = '  +1.50%  '
$ws.Range("E22").Value = '  +1.01%  '
$ws.Range("D23").Value = '''0.0000190'
$ws.Range("E23").Value = '  +48.73%  '
$ws.Range("D24").Value = '''5.30'
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("D25").Value = '''5.39'
$ws.Range("E25").Value = '  +5.84%  '
$ws.Range("D26").Value = '''12.23'
$ws.Range("E26").Value = '  -1.85%  '
$ws.Range("D27").Value = '3.513.27'
$ws.Range("E27").Value = '  +2.82%  '
$ws.Range("D28").Value = '''76.68'
$ws.Range("E28").Value = '  +1.27%  '
$ws.Range("E29").Value = '  +0.16%  '
$ws.Range("D30").Value = '''0.181'
$ws.Range("E30").Value = '  +4.57%  '
$ws.Range("D31").Value = '''1.00'
$ws.Range("E31").Value = '  +0.04%  '
$ws.Range("D32").Value = '''8.81'
$ws.Range("E32").Value = '  +1.27%  '
$ws.Range("D33").Value = '''560.77'
$ws.Range("E33").Value = '  +3.45%  '
$ws.Range("D34").Value = '''7.16'
$ws.Range("E34").Value = '  +2.09%  '
$ws.Range("E35").Value = '  +25.97%  '
$ws.Range("E36").Value = '  +0.05%  '
$ws.Range("E37").Value = '  -7.41%  '
$ws.Range("D38").Value = '''22.75'
$ws.Range("E38").Value = '  +1.98%  '
$ws.Range("E39").Value = '  -2.91%  '
$ws.Range("D40").Value = '''22.47'
$ws.Range("E40").Value = '  +3.87%  '
$ws.Range("D41").Value = '''1.00'
$ws.Range("E41").Value = '  +0.14%  '
$ws.Range("E42").Value = '  +1.57%  '
$ws.Range("D43").Value = '''2.00'
$ws.Range("E43").Value = '  +1.04%  '
$ws.Range("E44").Value = '  +0.13%  '
$ws.Range("D45").Value = '''152.15'
$ws.Range("E45").Value = '  -1.95%  '
$ws.Range("D46").Value = '''181.59'
$ws.Range("E46").Value = '  +1.56%  '
$ws.Range("D47").Value = '''43.92'
$ws.Range("E47").Value = '  -1.91%  '
$ws.Range("E48").Value = '  +5.70%  '
$ws.Range("E49").Value = '  +0.50%  '
$ws.Range("D50").Value = '''0.635'
$ws.Range("E51").Value = '  +0.19%  '
